$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Waves_009.txt"

$ws.Range("D2").Value = 72
$ws.Range("E2").Value = 17
$ws.Range("F2").Value = 449.43
$ws.Range("G2").Value = 7.54
$ws.Range("H2").Value = 4.9
$ws.Range("I2").Value = 3.62
$ws.Range("J2").Value = 0.1
$ws.Range("K2").Value = 3.38
$ws.Range("L2").Value = 2.1
$ws.Range("M2").Value = 0.1
$ws.Range("N2").Value = 2.25
$ws.Range("O2").Value = 0.82
$ws.Range("P2").Value = 0.04
$ws.Range("Q2").Value = 482.29
$ws.Range("R2").Value = 540.05
$ws.Range("S2").Value = 14.68
$ws.Range("T2").Value = 14.07
$ws.Range("U2").Value = 18.61
$ws.Range("V2").Value = 0.51
$ws.Range("W2").Value = 526.6
$ws.Range("X2").Value = 596.0700000000001
$ws.Range("Y2").Value = 16.2
$ws.Range("Z2").Value = 14.31
$ws.Range("AA2").Value = 16.74
$ws.Range("AB2").Value = 0.46
$ws.Range("AC2").Value = 3.17
$ws.Range("AD2").Value = 1.3
$ws.Range("AF2").Value = 8.93
$ws.Range("AG2").Value = 3.52
$ws.Range("AH2").Value = 0.1
$ws.Range("AI2").Value = 13.46
$ws.Range("AJ2").Value = 4.92
$ws.Range("AK2").Value = 0.13
